# Auto-generated script applying scheduled-runner price updates
# to the Atomos_Profits workbook (columns H..N per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Cells.Item(6, 8).Value2 = 1871.5
$ws.Cells.Item(6, 9).Value2 = 614.5
$ws.Cells.Item(6, 10).Value2 = 2500
$ws.Cells.Item(6, 11).Value2 = 1843.5
$ws.Cells.Item(6, 12).Value2 = 7500
$ws.Cells.Item(6, 13).Value2 = -1731.5
$ws.Cells.Item(6, 14).Value2 = -7724
# ALC row 11
$ws.Cells.Item(11, 8).Value2 = 1165.5
$ws.Cells.Item(11, 9).Value2 = 1165.5
$ws.Cells.Item(11, 11).Value2 = 1165.5
$ws.Cells.Item(11, 13).Value2 = -1025.5
# ALC row 31
$ws.Cells.Item(31, 8).Value2 = 1139.6
$ws.Cells.Item(31, 9).Value2 = 924.5
$ws.Cells.Item(31, 10).Value2 = 2000
$ws.Cells.Item(31, 11).Value2 = 2773.5
$ws.Cells.Item(31, 12).Value2 = 6000
$ws.Cells.Item(31, 13).Value2 = -2543.5
$ws.Cells.Item(31, 14).Value2 = -6460
# ALC row 38
$ws.Cells.Item(38, 8).Value2 = 990.5238000000001
$ws.Cells.Item(38, 9).Value2 = 220.2
$ws.Cells.Item(38, 10).Value2 = 2916.3333
$ws.Cells.Item(38, 11).Value2 = 660.5999999999999
$ws.Cells.Item(38, 12).Value2 = 8748.999899999999
$ws.Cells.Item(38, 13).Value2 = -288.5999999999999
$ws.Cells.Item(38, 14).Value2 = -9492.999899999999
# ALC row 39
$ws.Cells.Item(39, 8).Value2 = 650.9231
$ws.Cells.Item(39, 9).Value2 = 45.3
$ws.Cells.Item(39, 10).Value2 = 2669.6667
$ws.Cells.Item(39, 11).Value2 = 135.9
$ws.Cells.Item(39, 12).Value2 = 8009.000100000001
$ws.Cells.Item(39, 13).Value2 = 160.1
$ws.Cells.Item(39, 14).Value2 = -8601.000100000001
# ALC row 137
$ws.Cells.Item(137, 8).Value2 = 2780.0222
$ws.Cells.Item(137, 9).Value2 = 3288.261
$ws.Cells.Item(137, 10).Value2 = 2248.682
$ws.Cells.Item(137, 11).Value2 = 9864.782999999999
$ws.Cells.Item(137, 12).Value2 = 6746.045999999999
$ws.Cells.Item(137, 13).Value2 = -7314.782999999999
$ws.Cells.Item(137, 14).Value2 = -11846.046

$ws = $wb.Worksheets.Item("ARM")
# ARM row 10
$ws.Cells.Item(10, 8).Value2 = 60202.4
$ws.Cells.Item(10, 9).Value2 = 41004
$ws.Cells.Item(10, 10).Value2 = 65002
$ws.Cells.Item(10, 11).Value2 = 41004
$ws.Cells.Item(10, 12).Value2 = 65002
$ws.Cells.Item(10, 13).Value2 = -40834
$ws.Cells.Item(10, 14).Value2 = -65342
# ARM row 19
$ws.Cells.Item(19, 8).Value2 = 65006.75
$ws.Cells.Item(19, 9).Value2 = 20000
$ws.Cells.Item(19, 10).Value2 = 80009
$ws.Cells.Item(19, 11).Value2 = 20000
$ws.Cells.Item(19, 12).Value2 = 80009
$ws.Cells.Item(19, 13).Value2 = -19771
$ws.Cells.Item(19, 14).Value2 = -80467
# ARM row 36
$ws.Cells.Item(36, 8).Value2 = 31710.2
$ws.Cells.Item(36, 9).Value2 = 6164.3335
$ws.Cells.Item(36, 11).Value2 = 6164.3335
$ws.Cells.Item(36, 13).Value2 = -5818.3335
# ARM row 61
$ws.Cells.Item(61, 8).Value2 = 4213.4116
$ws.Cells.Item(61, 9).Value2 = 1519.5
$ws.Cells.Item(61, 10).Value2 = 5042.3076
$ws.Cells.Item(61, 11).Value2 = 1519.5
$ws.Cells.Item(61, 12).Value2 = 5042.3076
$ws.Cells.Item(61, 13).Value2 = -1307.5
$ws.Cells.Item(61, 14).Value2 = -5466.3076
# ARM row 63
$ws.Cells.Item(63, 8).Value2 = 3611.111
$ws.Cells.Item(63, 9).Value2 = 2500
$ws.Cells.Item(63, 10).Value2 = 4500
$ws.Cells.Item(63, 11).Value2 = 2500
$ws.Cells.Item(63, 12).Value2 = 4500
$ws.Cells.Item(63, 13).Value2 = -1814
$ws.Cells.Item(63, 14).Value2 = -5872
# ARM row 66
$ws.Cells.Item(66, 8).Value2 = 3611.111
$ws.Cells.Item(66, 9).Value2 = 2500
$ws.Cells.Item(66, 10).Value2 = 4500
$ws.Cells.Item(66, 11).Value2 = 12500
$ws.Cells.Item(66, 12).Value2 = 22500
$ws.Cells.Item(66, 13).Value2 = -9068
$ws.Cells.Item(66, 14).Value2 = -29364
# ARM row 103
$ws.Cells.Item(103, 8).Value2 = 25013.408
$ws.Cells.Item(103, 10).Value2 = 25013.408
$ws.Cells.Item(103, 12).Value2 = 25013.408
$ws.Cells.Item(103, 14).Value2 = -27357.408
# ARM row 136
$ws.Cells.Item(136, 8).Value2 = 4213.4116
$ws.Cells.Item(136, 9).Value2 = 1519.5
$ws.Cells.Item(136, 10).Value2 = 5042.3076
$ws.Cells.Item(136, 11).Value2 = 4558.5
$ws.Cells.Item(136, 12).Value2 = 15126.9228
$ws.Cells.Item(136, 13).Value2 = -2008.5
$ws.Cells.Item(136, 14).Value2 = -20226.9228

$ws = $wb.Worksheets.Item("CRP")
# CRP row 2
$ws.Cells.Item(2, 8).Value2 = 44703
$ws.Cells.Item(2, 9).Value2 = 1500
$ws.Cells.Item(2, 10).Value2 = 55503.75
$ws.Cells.Item(2, 11).Value2 = 1500
$ws.Cells.Item(2, 12).Value2 = 55503.75
$ws.Cells.Item(2, 13).Value2 = -1387
$ws.Cells.Item(2, 14).Value2 = -55729.75
# CRP row 11
$ws.Cells.Item(11, 8).Value2 = 69006
$ws.Cells.Item(11, 10).Value2 = 69006
$ws.Cells.Item(11, 12).Value2 = 69006
$ws.Cells.Item(11, 14).Value2 = -69286
# CRP row 17
$ws.Cells.Item(17, 8).Value2 = 46621
$ws.Cells.Item(17, 9).Value2 = 29954
$ws.Cells.Item(17, 10).Value2 = 54954.5
$ws.Cells.Item(17, 11).Value2 = 29954
$ws.Cells.Item(17, 12).Value2 = 54954.5
$ws.Cells.Item(17, 13).Value2 = -29780
$ws.Cells.Item(17, 14).Value2 = -55302.5
# CRP row 31
$ws.Cells.Item(31, 8).Value2 = 2702.1516
$ws.Cells.Item(31, 9).Value2 = 1998
$ws.Cells.Item(31, 11).Value2 = 1998
$ws.Cells.Item(31, 13).Value2 = -1703
# CRP row 32
$ws.Cells.Item(32, 8).Value2 = 46003.2
$ws.Cells.Item(32, 9).Value2 = 10005
$ws.Cells.Item(32, 10).Value2 = 55002.75
$ws.Cells.Item(32, 11).Value2 = 10005
$ws.Cells.Item(32, 12).Value2 = 55002.75
$ws.Cells.Item(32, 13).Value2 = -9689
$ws.Cells.Item(32, 14).Value2 = -55634.75
# CRP row 34
$ws.Cells.Item(34, 8).Value2 = 2702.1516
$ws.Cells.Item(34, 9).Value2 = 1998
$ws.Cells.Item(34, 11).Value2 = 1998
$ws.Cells.Item(34, 13).Value2 = -1796
# CRP row 45
$ws.Cells.Item(45, 8).Value2 = 10022.333
$ws.Cells.Item(45, 10).Value2 = 5000
$ws.Cells.Item(45, 12).Value2 = 5000
$ws.Cells.Item(45, 14).Value2 = -6186
# CRP row 105
$ws.Cells.Item(105, 8).Value2 = 2848.9048
$ws.Cells.Item(105, 9).Value2 = 2554.0527
$ws.Cells.Item(105, 11).Value2 = 2554.0527
$ws.Cells.Item(105, 13).Value2 = -807.0527000000002
# CRP row 134
$ws.Cells.Item(134, 8).Value2 = 3613.5334
$ws.Cells.Item(134, 9).Value2 = 2049
$ws.Cells.Item(134, 10).Value2 = 5401.5713
$ws.Cells.Item(134, 11).Value2 = 6147
$ws.Cells.Item(134, 12).Value2 = 16204.7139
$ws.Cells.Item(134, 13).Value2 = -3612
$ws.Cells.Item(134, 14).Value2 = -21274.7139

$ws = $wb.Worksheets.Item("CUL")
# CUL row 7
$ws.Cells.Item(7, 8).Value2 = 195
$ws.Cells.Item(7, 9).Value2 = 190
$ws.Cells.Item(7, 10).Value2 = 200
$ws.Cells.Item(7, 11).Value2 = 570
$ws.Cells.Item(7, 12).Value2 = 600
$ws.Cells.Item(7, 13).Value2 = -458
$ws.Cells.Item(7, 14).Value2 = -824
# CUL row 131
$ws.Cells.Item(131, 8).Value2 = 1431.1818
$ws.Cells.Item(131, 9).Value2 = 1744.375
$ws.Cells.Item(131, 10).Value2 = 1252.2142
$ws.Cells.Item(131, 11).Value2 = 5233.125
$ws.Cells.Item(131, 12).Value2 = 3756.6426
$ws.Cells.Item(131, 13).Value2 = -193.125
$ws.Cells.Item(131, 14).Value2 = -13836.6426

$ws = $wb.Worksheets.Item("GSM")
# GSM row 102
$ws.Cells.Item(102, 8).Value2 = 2491.3
$ws.Cells.Item(102, 9).Value2 = 1678.9231
$ws.Cells.Item(102, 11).Value2 = 1678.9231
$ws.Cells.Item(102, 13).Value2 = -56.92309999999998

$ws = $wb.Worksheets.Item("LTW")
# LTW row 32
$ws.Cells.Item(32, 8).Value2 = 29132.625
$ws.Cells.Item(32, 9).Value2 = 3254
$ws.Cells.Item(32, 10).Value2 = 55011.25
$ws.Cells.Item(32, 11).Value2 = 3254
$ws.Cells.Item(32, 12).Value2 = 55011.25
$ws.Cells.Item(32, 13).Value2 = -2937
$ws.Cells.Item(32, 14).Value2 = -55645.25
# LTW row 34
$ws.Cells.Item(34, 8).Value2 = 0
$ws.Cells.Item(34, 10).Value2 = 0
$ws.Cells.Item(34, 12).Value2 = 0
$ws.Cells.Item(34, 14).ClearContents() | Out-Null
# LTW row 136
$ws.Cells.Item(136, 8).Value2 = 1912.8214
$ws.Cells.Item(136, 9).Value2 = 1506.0555
$ws.Cells.Item(136, 10).Value2 = 2645
$ws.Cells.Item(136, 11).Value2 = 4518.166499999999
$ws.Cells.Item(136, 12).Value2 = 7935
$ws.Cells.Item(136, 13).Value2 = -1968.166499999999
$ws.Cells.Item(136, 14).Value2 = -13035

$ws = $wb.Worksheets.Item("WVR")
# WVR row 13
$ws.Cells.Item(13, 8).Value2 = 37198.4
$ws.Cells.Item(13, 10).Value2 = 37198.4
$ws.Cells.Item(13, 12).Value2 = 37198.4
$ws.Cells.Item(13, 14).Value2 = -37478.4
# WVR row 17
$ws.Cells.Item(17, 8).Value2 = 0
$ws.Cells.Item(17, 9).Value2 = 0
$ws.Cells.Item(17, 10).Value2 = 0
$ws.Cells.Item(17, 11).Value2 = 0
$ws.Cells.Item(17, 12).Value2 = 0
$ws.Cells.Item(17, 13).ClearContents() | Out-Null
$ws.Cells.Item(17, 14).ClearContents() | Out-Null
# WVR row 23
$ws.Cells.Item(23, 8).Value2 = 48396.6
$ws.Cells.Item(23, 9).Value2 = 975
$ws.Cells.Item(23, 10).Value2 = 80011
$ws.Cells.Item(23, 11).Value2 = 975
$ws.Cells.Item(23, 12).Value2 = 80011
$ws.Cells.Item(23, 13).Value2 = -746
$ws.Cells.Item(23, 14).Value2 = -80469
# WVR row 41
$ws.Cells.Item(41, 8).Value2 = 6867.5
$ws.Cells.Item(41, 9).Value2 = 5342
$ws.Cells.Item(41, 11).Value2 = 5342
$ws.Cells.Item(41, 13).Value2 = -4952
# WVR row 58
$ws.Cells.Item(58, 8).Value2 = 14010.625
$ws.Cells.Item(58, 10).Value2 = 14000
$ws.Cells.Item(58, 12).Value2 = 14000
$ws.Cells.Item(58, 14).Value2 = -14616

